$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.206.51'
$ws.Range('D3').Value = '3.401.06'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.52'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.87'
$ws.Range('E6').Value = '  +1.80%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.400.34'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.69'
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.380'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '3.982.74'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.72'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D17').Value = '3.406.16'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '61.254.67'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.93'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.35'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.63'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').Value = '3.530.47'
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.552'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.12'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('E28').Value = '  -4.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.173'
$ws.Range('E29').Value = '  +8.18%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.11'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.43'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '166.22'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0769'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.07'
$ws.Range('E41').Value = '  +6.26%  '
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.778'
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.95'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.38'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').Value = '2.522.89'
$ws.Range('E48').Value = '  +7.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.70'
$ws.Range('E49').Value = '  +4.93%  '
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0261'
$ws.Range('E51').Value = '  +0.13%  '
